$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column A throughout this log uses the short-date display (style index 1,
# numFmtId 14 "m/d/yy" in this workbook) - reapply it to every A cell we
# touch below so new rows match the rest of the column.
$dateFmt = "m/d/yy"

# --- Row 176: becomes the former row-177 "Sine3" data row (shifted up) ---
$ws.Range("A176").Value = 44650
$ws.Range("A176").NumberFormat = $dateFmt
$ws.Range("B176").Value = "Sine3"
$ws.Range("C176").Value = 6
$ws.Range("D176").Value = 3
$ws.Range("E176").Value = 6
$ws.Range("F176").Value = 2

# --- Rows 177-181: new Ramps3 trials (Trial 3..7) ---
$ws.Range("A177").Value = 44650
$ws.Range("A177").NumberFormat = $dateFmt
$ws.Range("B177").Value = "Ramps3"
$ws.Range("C177").Value = 3
$ws.Range("D177").Value = 4
$ws.Range("E177").Value = 4
$ws.Range("F177").Value = ""

$ws.Range("A178").Value = 44650
$ws.Range("A178").NumberFormat = $dateFmt
$ws.Range("B178").Value = "Ramps3"
$ws.Range("C178").Value = 4
$ws.Range("D178").Value = 4
$ws.Range("E178").Value = 4

$ws.Range("A179").Value = 44650
$ws.Range("A179").NumberFormat = $dateFmt
$ws.Range("B179").Value = "Ramps3"
$ws.Range("C179").Value = 5
$ws.Range("D179").Value = 4
$ws.Range("E179").Value = 4

$ws.Range("A180").Value = 44650
$ws.Range("A180").NumberFormat = $dateFmt
$ws.Range("B180").Value = "Ramps3"
$ws.Range("C180").Value = 6
$ws.Range("D180").Value = 4
$ws.Range("E180").Value = 4

$ws.Range("A181").Value = 44650
$ws.Range("A181").NumberFormat = $dateFmt
$ws.Range("B181").Value = "Ramps3"
$ws.Range("C181").Value = 7
$ws.Range("D181").Value = 4
$ws.Range("E181").Value = 4

# --- Row 182: note (reuses the old "Rewired load cells..." string slot, retexted) ---
$ws.Range("A182").Value = "Changed sign on drawWire"
$ws.Range("A182").NumberFormat = $dateFmt

# --- Row 183: Ramps3 Trial 8 ---
$ws.Range("A183").Value = 44650
$ws.Range("A183").NumberFormat = $dateFmt
$ws.Range("B183").Value = "Ramps3"
$ws.Range("C183").Value = 8
$ws.Range("D183").Value = 4
$ws.Range("E183").Value = 4

# --- Row 184: new note ---
$ws.Range("A184").Value = "Added offset to drawWire"
$ws.Range("A184").NumberFormat = $dateFmt

# --- Row 185: Ramps3 Trial 9 ---
$ws.Range("A185").Value = 44650
$ws.Range("A185").NumberFormat = $dateFmt
$ws.Range("B185").Value = "Ramps3"
$ws.Range("C185").Value = 9
$ws.Range("D185").Value = 4
$ws.Range("E185").Value = 4

# --- Rows 186-188: Ramps3 Trials 10-12, dated 44655 ---
$ws.Range("A186").Value = 44655
$ws.Range("A186").NumberFormat = $dateFmt
$ws.Range("B186").Value = "Ramps3"
$ws.Range("C186").Value = 10
$ws.Range("D186").Value = 4
$ws.Range("E186").Value = 4

$ws.Range("A187").Value = 44655
$ws.Range("A187").NumberFormat = $dateFmt
$ws.Range("B187").Value = "Ramps3"
$ws.Range("C187").Value = 11
$ws.Range("D187").Value = 4
$ws.Range("E187").Value = 4

$ws.Range("A188").Value = 44655
$ws.Range("A188").NumberFormat = $dateFmt
$ws.Range("B188").Value = "Ramps3"
$ws.Range("C188").Value = 12
$ws.Range("D188").Value = 4
$ws.Range("E188").Value = 4

# --- Remove the last row (448) entirely; the used range shrinks to A1:H447 ---
$ws.Range("A448").EntireRow.Delete()

# --- Restore the selection recorded in the file ---
$ws.Range("G188").Select()
